$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the table one column to the left: B,C,D -> A,B,C
# and correct the "Edad" header typo to "Edada1" at the same time.
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Apellido"
$ws.Range("C1").Value = "Edada1"

$ws.Range("A2").Value = "yesid"
$ws.Range("B2").Value = "ochoa"
$ws.Range("C2").Value = 25

$ws.Range("A3").Value = "valentina"
$ws.Range("B3").Value = "ramos"
$ws.Range("C3").Value = 26

$ws.Range("A4").Value = "carlos"
$ws.Range("B4").Value = "luque"
$ws.Range("C4").Value = 27

$ws.Range("A5").Value = "Nataly"
$ws.Range("B5").Value = "luque"
$ws.Range("C5").Value = 28

# The old column D, now superseded by the shifted column C, is cleared.
$ws.Range("D1:D5").ClearContents()

# Add the new underlined (blank) cell block in I14:L18.
$ws.Range("I14:L18").Font.Underline = $true

# Move the active selection to F3, matching the edited sheet view.
[void]$ws.Range("F3").Select()
